$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text so values like "0.520"
# or "215.51" are not auto-converted to numbers and keep their exact
# original formatting (trailing zeros, thousands dot-separators, etc.)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.067.33'
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.657.49'
$ws.Range("E3").Value = '  +3.09%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.51'
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.513'
$ws.Range("E6").Value = '  +2.30%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").Value = '  +2.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0616'
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.28'
$ws.Range("E10").Value = '  +4.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0876'
$ws.Range("E11").Value = '  +2.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.890.15'
$ws.Range("E12").Value = '  +3.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.661.65'
$ws.Range("E13").Value = '  +3.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.08'
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.520'
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.37'
$ws.Range("E16").Value = '  +2.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.046.00'
$ws.Range("E17").Value = '  +2.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '237.39'
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.76'
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0732'
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.44'
$ws.Range("E22").Value = '  +3.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.34'
$ws.Range("E23").Value = '  +3.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.21'
$ws.Range("E24").Value = '  +3.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.63'
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.12'
$ws.Range("E26").Value = '  +1.54%  '
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.88'
$ws.Range("E28").Value = '  +2.37%  '
$ws.Range("B29").Value = 'BinanceUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("E31").Value = '  +1.64%  '
$ws.Range("B32").Value = 'Maker'
$ws.Range("C32").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.571.40'
$ws.Range("E32").Value = '  +4.56%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.33'
$ws.Range("E33").Value = '  +2.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.09'
$ws.Range("E34").Value = '  +4.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").Value = '  +7.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.41'
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.583'
$ws.Range("E37").Value = '  +2.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.906'
$ws.Range("E38").Value = '  +9.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0170'
$ws.Range("E39").Value = '  +2.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.98'
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("E42").Value = '  +3.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.61'
$ws.Range("E43").Value = '  +7.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.797.18'
$ws.Range("E44").Value = '  +3.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.777'
$ws.Range("E45").Value = '  +1.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.912'
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.27'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.52'
$ws.Range("E49").Value = '  +1.20%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0986'
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0506'
$ws.Range("E51").Value = '  +1.12%  '
